# Apply the "Making the code easier to read" edit to burnrate-simple.xlsx
# Substantive changes only (cosmetic app/version metadata left untouched):
#   - Active/selected sheet moves from "Intro" to "Propellant"
#   - New selections recorded on Propellant / a-n calcs / C-Star
#   - ProPep view scrolled down (topLeftCell = A10)
#   - Propellant!C4 becomes a formula (=0.9*C5) instead of a hard-coded value
#   - Several input cells change value across Propellant, a-n calcs and C-Star

$wb = $excel.ActiveWorkbook

$wsIntro      = $wb.Worksheets.Item("Intro")
$wsPropellant = $wb.Worksheets.Item("Propellant")
$wsAN         = $wb.Worksheets.Item("a-n calcs")
$wsCStar      = $wb.Worksheets.Item("C-Star")
$wsProPep     = $wb.Worksheets.Item("ProPep")

# --- Propellant sheet: update input values -------------------------------
$wsPropellant.Range("C5").Value = 1.8828
$wsPropellant.Range("C4").Formula = "=0.9*C5"
$wsPropellant.Range("C11").Value = 1720
$wsPropellant.Range("C13").Value = 2546
$wsPropellant.Range("C15").Value = 1.133

# --- a-n calcs sheet: update input value -----------------------------------
$wsAN.Range("E6").Value = 1.6

# --- C-Star sheet: clear some inputs, fill in others -----------------------
$wsCStar.Range("C5").Value = $null
$wsCStar.Range("C6").Value = 21.659
$wsCStar.Range("C10").Value = $null
$wsCStar.Range("C11").Value = 0.34468
$wsCStar.Range("C13").Value = $null
$wsCStar.Range("C14").Value = 12.54
$wsCStar.Range("C16").Value = 0.1

# --- View / selection state -------------------------------------------------
$wsProPep.Activate()
$wsProPep.Application.ActiveWindow.ScrollRow = 10

$wsAN.Activate()
$wsAN.Range("E6").Select()

$wsCStar.Activate()
$wsCStar.Range("C16").Select()

$wsPropellant.Activate()
$wsPropellant.Range("C15").Select()

$wb.Activate()
